# Finished re-running SigProfilerExtractor for SBS_set2/Realistic on HPC-cluster.
# The newly-generated results (the "seed.1076753" run) are inserted into the
# SBS_set2 / SigProfilerExtractor block of the summary table, pushing the
# later rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The SBS_set2 / SigProfilerExtractor block currently starts at row 77
# (seed.145879, seed.200437, seed.310111, seed.528401). Insert a new row
# above it for the newly completed "seed.1076753" run, shifting every row
# from 77 downward (through the previous last row, 90) down by one.
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row with the freshly generated result.
$ws.Cells.Item(77, 1).Value = "SBS_set2"
$ws.Cells.Item(77, 2).Value = "SigProfilerExtractor"
$ws.Cells.Item(77, 3).Value = "seed.1076753"
$ws.Cells.Item(77, 4).Value = 2102332.48
